$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Will"
$ws.Range("A2").Value = "is"
$ws.Range("A3").Value = "cool"
